$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update precipitation values for B2:B7
$ws.Range("B2:B7").Value = 445.00799999999998

# Move the active selection/cursor
$ws.Range("E13").Select()

# Set page orientation to portrait (xlPortrait = 1)
$ws.PageSetup.Orientation = 1
